# Remove the LTE antenna (ANT2 / FXUB65.07.0180C) line item from the BOM.
# This was row 2 of the sheet; deleting it shifts every subsequent row up
# by one, which also requires updating the dependent dimension/name/view
# state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the antenna row entirely (cells shift up).
$ws.Rows(2).Delete()

# Match the post-edit selection recorded in the workbook (row 2 selected).
$ws.Rows(2).EntireRow.Select() | Out-Null

# The named range "permamote_lte" covered the whole table; shrink it by
# one row now that the sheet is one row shorter.
$n = $wb.Names.Item(1)
$n.RefersTo = "=Sheet1!`$A`$1:`$S`$52"
